$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (new date, new EBITDA value). EBITDA is $null when unchanged.
$updates = @(
    @{ Row = 2;  Date = "2026/01/03"; Ebitda = "7.05"  },
    @{ Row = 8;  Date = "2026/01/03"; Ebitda = "8.55"  },
    @{ Row = 14; Date = "2026/01/03"; Ebitda = "3.07"  },
    @{ Row = 20; Date = "2026/01/03"; Ebitda = "12.97" },
    @{ Row = 26; Date = "2026/01/03"; Ebitda = "11.47" },
    @{ Row = 32; Date = "2026/01/03"; Ebitda = "27.59" },
    @{ Row = 38; Date = "2026/01/03"; Ebitda = $null   },
    @{ Row = 44; Date = "2026/01/03"; Ebitda = "11.63" },
    @{ Row = 50; Date = "2026/01/03"; Ebitda = "11.35" },
    @{ Row = 56; Date = "2026/01/03"; Ebitda = "31.65" },
    @{ Row = 62; Date = "2026/01/03"; Ebitda = "11.96" },
    @{ Row = 68; Date = "2026/01/03"; Ebitda = "13.24" },
    @{ Row = 74; Date = "2026/01/03"; Ebitda = "17.07" }
)

foreach ($u in $updates) {
    $r = $u.Row

    # Column A holds the date-like string - force it to stay plain text
    # (not an Excel date serial) the same way the source data was stored.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $u.Date
    $cellA.ClearFormats()

    if ($u.Ebitda -ne $null) {
        # Column B holds a numeric-looking string - force it to stay text too.
        $cellB = $ws.Cells.Item($r, 2)
        $cellB.NumberFormat = "@"
        $cellB.Value = $u.Ebitda
        $cellB.ClearFormats()
    }
}
